$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying records for three pairs of rows were re-sorted; the row
# numbers/positions stay fixed but the record content (cell values) swaps
# between each pair. Swap cell-by-cell for the columns that actually differ.

$pairs = @(
    @{ Row1 = 13; Row2 = 15; Cols = @("A","B","E","F","G","H","Q","R","Z","AB") },
    @{ Row1 = 22; Row2 = 23; Cols = @("A","B","D","E","F","G","H","Q","R","S","Z","AB","AC") },
    @{ Row1 = 27; Row2 = 28; Cols = @("A","B","E","F","G","H","Q","R","S","Z","AB","AC") }
)

foreach ($pair in $pairs) {
    $row1 = $pair.Row1
    $row2 = $pair.Row2
    $cols = $pair.Cols
    foreach ($col in $cols) {
        $addr1 = $col + $row1
        $addr2 = $col + $row2
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}
